$wb = $excel.ActiveWorkbook

# ---- ALC (Worksheets.Item(1)) ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("H17").Value = 663.8617
$ws.Range("J17").Value = 668.8495
$ws.Range("L17").Value = 2006.5485
$ws.Range("N17").Value = -2342.5485

$ws.Range("H33").Value = 1237307.9
$ws.Range("J33").Value = 1500
$ws.Range("L33").Value = 1500
$ws.Range("N33").Value = -1958

$ws.Range("H129").Value = 670.1905
$ws.Range("I129").Value = 382.8421
$ws.Range("J129").Value = 3400
$ws.Range("K129").Value = 1148.5263
$ws.Range("L129").Value = 10200
$ws.Range("M129").Value = 3851.4737
$ws.Range("N129").Value = -20200

$ws.Range("H132").Value = 1932.7084
$ws.Range("I132").Value = 1276.7894
$ws.Range("J132").Value = 4425.2
$ws.Range("K132").Value = 3830.3682
$ws.Range("L132").Value = 13275.6
$ws.Range("M132").Value = -1300.3682
$ws.Range("N132").Value = -18335.6

$ws.Range("H135").Value = 35664.07
$ws.Range("I135").Value = 44486.957
$ws.Range("K135").Value = 400382.613
$ws.Range("M135").Value = -397847.613

$ws.Range("H137").Value = 3262031.8
$ws.Range("I137").Value = 1516308
$ws.Range("J137").Value = 7693484.5
$ws.Range("K137").Value = 4548924
$ws.Range("L137").Value = 23080453.5
$ws.Range("M137").Value = -4546374
$ws.Range("N137").Value = -23085553.5


# ---- ARM (Worksheets.Item(2)) ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 2648819
$ws.Range("I32").Value = 4225.4937
$ws.Range("J32").Value = 17949680
$ws.Range("K32").Value = 4225.4937
$ws.Range("L32").Value = 17949680
$ws.Range("M32").Value = -3938.4937
$ws.Range("N32").Value = -17950254

$ws.Range("H74").Value = 873.0208
$ws.Range("I74").Value = 880.5349
$ws.Range("J74").Value = 808.4
$ws.Range("K74").Value = 880.5349
$ws.Range("L74").Value = 808.4
$ws.Range("M74").Value = -6.534899999999993
$ws.Range("N74").Value = -2556.4

$ws.Range("H77").Value = 873.0208
$ws.Range("I77").Value = 880.5349
$ws.Range("J77").Value = 808.4
$ws.Range("K77").Value = 4402.6745
$ws.Range("L77").Value = 4042
$ws.Range("M77").Value = -34.67450000000008
$ws.Range("N77").Value = -12778

$ws.Range("H113").Value = 49195
$ws.Range("J113").Value = 49195
$ws.Range("L113").Value = 49195
$ws.Range("N113").Value = -57873

$ws.Range("H122").Value = 1511.5
$ws.Range("I122").Value = 1481.9656
$ws.Range("J122").Value = 1682.8
$ws.Range("K122").Value = 4445.8968
$ws.Range("L122").Value = 5048.4
$ws.Range("M122").Value = -1995.8968
$ws.Range("N122").Value = -9948.4

$ws.Range("H132").Value = 2148.9424
$ws.Range("I132").Value = 1432.175
$ws.Range("J132").Value = 4538.1665
$ws.Range("K132").Value = 4296.525
$ws.Range("L132").Value = 13614.4995
$ws.Range("M132").Value = -1766.525
$ws.Range("N132").Value = -18674.4995


# ---- CRP (Worksheets.Item(4)) ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 2872
$ws.Range("I31").Value = 2086.923
$ws.Range("J31").Value = 5423.5
$ws.Range("K31").Value = 2086.923
$ws.Range("L31").Value = 5423.5
$ws.Range("M31").Value = -1791.923
$ws.Range("N31").Value = -6013.5

$ws.Range("H34").Value = 2872
$ws.Range("I34").Value = 2086.923
$ws.Range("J34").Value = 5423.5
$ws.Range("K34").Value = 2086.923
$ws.Range("L34").Value = 5423.5
$ws.Range("M34").Value = -1884.923
$ws.Range("N34").Value = -5827.5

$ws.Range("H63").Value = 33000
$ws.Range("J63").Value = 33000
$ws.Range("L63").Value = 33000
$ws.Range("N63").Value = -34372

$ws.Range("H66").Value = 33000
$ws.Range("J66").Value = 33000
$ws.Range("L66").Value = 99000
$ws.Range("N66").Value = -105864

$ws.Range("H81").Value = 32000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 32000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 32000
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -33996

$ws.Range("H84").Value = 32000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 32000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 96000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -105984

$ws.Range("H118").Value = 28000
$ws.Range("J118").Value = 28000
$ws.Range("L118").Value = 28000
$ws.Range("N118").Value = -31314

$ws.Range("H134").Value = 2889.697
$ws.Range("I134").Value = 2903.4333
$ws.Range("J134").Value = 2752.3333
$ws.Range("K134").Value = 8710.2999
$ws.Range("L134").Value = 8256.999899999999
$ws.Range("M134").Value = -6175.2999
$ws.Range("N134").Value = -13326.9999


# ---- CUL (Worksheets.Item(5)) ----
$ws = $wb.Worksheets.Item(5)
$ws.Range("H117").Value = 5863.636
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 5863.636
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 17590.908
$ws.Range("M117").ClearContents()
$ws.Range("N117").Value = -24474.908

$ws.Range("H120").Value = 17268.182
$ws.Range("I120").Value = 4975
$ws.Range("J120").Value = 20000
$ws.Range("K120").Value = 14925
$ws.Range("L120").Value = 60000
$ws.Range("M120").Value = -10087
$ws.Range("N120").Value = -69676

$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()


# ---- GSM (Worksheets.Item(6)) ----
$ws = $wb.Worksheets.Item(6)
$ws.Range("H122").Value = 1757.4
$ws.Range("I122").Value = 1746.0625
$ws.Range("J122").Value = 1777.5555
$ws.Range("K122").Value = 5238.1875
$ws.Range("L122").Value = 5332.666499999999
$ws.Range("M122").Value = -2788.1875
$ws.Range("N122").Value = -10232.6665


# ---- LTW (Worksheets.Item(7)) ----
$ws = $wb.Worksheets.Item(7)
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H132").Value = 1884.804
$ws.Range("I132").Value = 1637.075
$ws.Range("J132").Value = 2785.6365
$ws.Range("K132").Value = 4911.225
$ws.Range("L132").Value = 8356.9095
$ws.Range("M132").Value = -2381.225
$ws.Range("N132").Value = -13416.9095


# ---- WVR (Worksheets.Item(8)) ----
$ws = $wb.Worksheets.Item(8)
$ws.Range("H46").Value = 49750
$ws.Range("J46").Value = 49750
$ws.Range("L46").Value = 49750
$ws.Range("N46").Value = -50212

$ws.Range("H132").Value = 1492.4219
$ws.Range("I132").Value = 1371.7322
$ws.Range("J132").Value = 2337.25
$ws.Range("K132").Value = 4115.196599999999
$ws.Range("L132").Value = 7011.75
$ws.Range("M132").Value = -1585.196599999999
$ws.Range("N132").Value = -12071.75

$ws.Range("H134").Value = 49750
$ws.Range("J134").Value = 49750
$ws.Range("L134").Value = 149250
$ws.Range("N134").Value = -154320

$ws.Range("H136").Value = 1130.625
$ws.Range("I136").Value = 1171.1014
$ws.Range("J136").Value = 199.66667
$ws.Range("K136").Value = 3513.3042
$ws.Range("L136").Value = 599.00001
$ws.Range("M136").Value = -963.3042
$ws.Range("N136").Value = -5699.00001

